$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10:F11").NumberFormat = "@"

$ws.Range("A10").Value = "Shriparna Gupta"
$ws.Range("B10").Value = "2001shriparna@gmail.com"
$ws.Range("C10").Value = "8420845844"
$ws.Range("D10").Value = "2020"
$ws.Range("E10").Value = "2969cc7d-a02b-426b-88d7-a46885ced627"
$ws.Range("F10").Value = ""

$ws.Range("A11").Value = "AGNIVA BHATTACHARJEE"
$ws.Range("B11").Value = "imagniva007@gmail.com"
$ws.Range("C11").Value = "08420880979"
$ws.Range("D11").Value = "2020"
$ws.Range("E11").Value = "c614ae38-316f-4b39-a47a-d16587c2f533"
$ws.Range("F11").Value = ""

$ws.Range("A1:F11").Errors.Item(9).Ignore = $true
